$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 2850.375
$ws.Cells.Item(40, 9).Value = 901
$ws.Cells.Item(40, 11).Value = 901
$ws.Cells.Item(40, 13).Value = -726

$ws.Cells.Item(41, 8).Value = 495.325
$ws.Cells.Item(41, 9).Value = 515.7619
$ws.Cells.Item(41, 10).Value = 472.73685
$ws.Cells.Item(41, 11).Value = 515.7619
$ws.Cells.Item(41, 12).Value = 472.73685
$ws.Cells.Item(41, 13).Value = -75.76189999999997
$ws.Cells.Item(41, 14).Value = -1352.73685

$ws.Cells.Item(112, 8).Value = 2181
$ws.Cells.Item(112, 10).Value = 2988
$ws.Cells.Item(112, 12).Value = 8964
$ws.Cells.Item(112, 14).Value = -11180

$ws.Cells.Item(113, 8).Value = 4454.9165
$ws.Cells.Item(113, 9).Value = 3291.8
$ws.Cells.Item(113, 11).Value = 3291.8
$ws.Cells.Item(113, 13).Value = -37.80000000000018

$ws.Cells.Item(132, 8).Value = 928584.4
$ws.Cells.Item(132, 9).Value = 1140552.9
$ws.Cells.Item(132, 10).Value = 10054.556
$ws.Cells.Item(132, 11).Value = 3421658.7
$ws.Cells.Item(132, 12).Value = 30163.668
$ws.Cells.Item(132, 13).Value = -3419128.7
$ws.Cells.Item(132, 14).Value = -35223.66800000001

$ws.Cells.Item(137, 8).Value = 54558.58
$ws.Cells.Item(137, 9).Value = 1781.9375
$ws.Cells.Item(137, 10).Value = 336034
$ws.Cells.Item(137, 11).Value = 5345.8125
$ws.Cells.Item(137, 12).Value = 1008102
$ws.Cells.Item(137, 13).Value = -2795.8125
$ws.Cells.Item(137, 14).Value = -1013202

$ws.Cells.Item(138, 8).Value = 4932.4346
$ws.Cells.Item(138, 10).Value = 5119.706
$ws.Cells.Item(138, 12).Value = 15359.118
$ws.Cells.Item(138, 14).Value = -25639.118

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2266.75
$ws.Cells.Item(2, 9).Value = 2327.3
$ws.Cells.Item(2, 10).Value = 2165.8333
$ws.Cells.Item(2, 11).Value = 2327.3
$ws.Cells.Item(2, 12).Value = 2165.8333
$ws.Cells.Item(2, 13).Value = -2214.3
$ws.Cells.Item(2, 14).Value = -2391.8333

$ws.Cells.Item(31, 8).Value = 5623.909
$ws.Cells.Item(31, 9).Value = 3134.4
$ws.Cells.Item(31, 10).Value = 30519
$ws.Cells.Item(31, 11).Value = 3134.4
$ws.Cells.Item(31, 12).Value = 30519
$ws.Cells.Item(31, 13).Value = -2840.4
$ws.Cells.Item(31, 14).Value = -31107

$ws.Cells.Item(32, 8).Value = 124289.9
$ws.Cells.Item(32, 9).Value = 123780.51
$ws.Cells.Item(32, 11).Value = 123780.51
$ws.Cells.Item(32, 13).Value = -123493.51

$ws.Cells.Item(45, 8).Value = 19705.055
$ws.Cells.Item(45, 9).Value = 18975.834
$ws.Cells.Item(45, 11).Value = 18975.834
$ws.Cells.Item(45, 13).Value = -18598.834

$ws.Cells.Item(97, 8).Value = 918.7368
$ws.Cells.Item(97, 9).Value = 814.375
$ws.Cells.Item(97, 11).Value = 814.375
$ws.Cells.Item(97, 13).Value = -318.375

$ws.Cells.Item(116, 8).Value = 2266.75
$ws.Cells.Item(116, 9).Value = 2327.3
$ws.Cells.Item(116, 10).Value = 2165.8333
$ws.Cells.Item(116, 11).Value = 2327.3
$ws.Cells.Item(116, 12).Value = 2165.8333
$ws.Cells.Item(116, 13).Value = -33.30000000000018
$ws.Cells.Item(116, 14).Value = -6753.8333

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2266.75
$ws.Cells.Item(3, 9).Value = 2327.3
$ws.Cells.Item(3, 10).Value = 2165.8333
$ws.Cells.Item(3, 11).Value = 2327.3
$ws.Cells.Item(3, 12).Value = 2165.8333
$ws.Cells.Item(3, 13).Value = -2213.3
$ws.Cells.Item(3, 14).Value = -2393.8333

$ws.Cells.Item(6, 8).Value = 80000
$ws.Cells.Item(6, 10).Value = 80000
$ws.Cells.Item(6, 12).Value = 80000
$ws.Cells.Item(6, 14).Value = -80226

$ws.Cells.Item(86, 8).Value = 1131.9756
$ws.Cells.Item(86, 9).Value = 1131.3513
$ws.Cells.Item(86, 11).Value = 1131.3513
$ws.Cells.Item(86, 13).Value = -8.351300000000037

$ws.Cells.Item(89, 8).Value = 1131.9756
$ws.Cells.Item(89, 9).Value = 1131.3513
$ws.Cells.Item(89, 11).Value = 5656.7565
$ws.Cells.Item(89, 13).Value = -40.75650000000041

$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 11).Value = 0
$ws.Cells.Item(93, 13).ClearContents()

$ws.Cells.Item(97, 8).Value = 6825.625
$ws.Cells.Item(97, 10).Value = 3993
$ws.Cells.Item(97, 12).Value = 3993
$ws.Cells.Item(97, 14).Value = -5975

$ws.Cells.Item(101, 8).Value = 0
$ws.Cells.Item(101, 10).Value = 0
$ws.Cells.Item(101, 12).Value = 0
$ws.Cells.Item(101, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2023.25
$ws.Cells.Item(16, 9).Value = 2204
$ws.Cells.Item(16, 10).Value = 1661.75
$ws.Cells.Item(16, 11).Value = 2204
$ws.Cells.Item(16, 12).Value = 1661.75
$ws.Cells.Item(16, 13).Value = -1917
$ws.Cells.Item(16, 14).Value = -2235.75

$ws.Cells.Item(31, 8).Value = 4220.3335
$ws.Cells.Item(31, 9).Value = 4220.3335
$ws.Cells.Item(31, 11).Value = 4220.3335
$ws.Cells.Item(31, 13).Value = -3925.3335

$ws.Cells.Item(34, 8).Value = 4220.3335
$ws.Cells.Item(34, 9).Value = 4220.3335
$ws.Cells.Item(34, 11).Value = 4220.3335
$ws.Cells.Item(34, 13).Value = -4018.3335

$ws.Cells.Item(60, 8).Value = 11009.1
$ws.Cells.Item(60, 9).Value = 10093
$ws.Cells.Item(60, 11).Value = 10093
$ws.Cells.Item(60, 13).Value = -9582

$ws.Cells.Item(97, 8).Value = 0
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 14).ClearContents()

$ws.Cells.Item(99, 8).Value = 9499.333000000001
$ws.Cells.Item(99, 9).Value = 8999
$ws.Cells.Item(99, 11).Value = 8999
$ws.Cells.Item(99, 13).Value = -7501

$ws.Cells.Item(113, 8).Value = 2023.25
$ws.Cells.Item(113, 9).Value = 2204
$ws.Cells.Item(113, 10).Value = 1661.75
$ws.Cells.Item(113, 11).Value = 2204
$ws.Cells.Item(113, 12).Value = 1661.75
$ws.Cells.Item(113, 13).Value = -34
$ws.Cells.Item(113, 14).Value = -6001.75

$ws.Cells.Item(126, 8).Value = 9499.333000000001
$ws.Cells.Item(126, 9).Value = 8999
$ws.Cells.Item(126, 11).Value = 26997
$ws.Cells.Item(126, 13).Value = -24527

$ws.Cells.Item(132, 8).Value = 4668.1304
$ws.Cells.Item(132, 9).Value = 4668.1304
$ws.Cells.Item(132, 11).Value = 14004.3912
$ws.Cells.Item(132, 13).Value = -11474.3912

$ws.Cells.Item(134, 8).Value = 2942.3333
$ws.Cells.Item(134, 9).Value = 2799
$ws.Cells.Item(134, 10).Value = 3014
$ws.Cells.Item(134, 11).Value = 8397
$ws.Cells.Item(134, 12).Value = 9042
$ws.Cells.Item(134, 13).Value = -5862
$ws.Cells.Item(134, 14).Value = -14112

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 190.66667
$ws.Cells.Item(12, 9).Value = 33.666668
$ws.Cells.Item(12, 10).Value = 347.66666
$ws.Cells.Item(12, 11).Value = 101.000004
$ws.Cells.Item(12, 12).Value = 1042.99998
$ws.Cells.Item(12, 13).Value = 71.999996
$ws.Cells.Item(12, 14).Value = -1388.99998

$ws.Cells.Item(33, 8).Value = 170.53334
$ws.Cells.Item(33, 9).Value = 60.333332
$ws.Cells.Item(33, 11).Value = 361.999992
$ws.Cells.Item(33, 13).Value = -78.99999200000002

$ws.Cells.Item(44, 8).Value = 2287.75
$ws.Cells.Item(44, 9).Value = 597
$ws.Cells.Item(44, 10).Value = 3198.1538
$ws.Cells.Item(44, 11).Value = 1791
$ws.Cells.Item(44, 12).Value = 9594.4614
$ws.Cells.Item(44, 13).Value = -1393
$ws.Cells.Item(44, 14).Value = -10390.4614

$ws.Cells.Item(64, 8).Value = 3228.5
$ws.Cells.Item(64, 10).Value = 3228.5
$ws.Cells.Item(64, 12).Value = 9685.5
$ws.Cells.Item(64, 14).Value = -10225.5

$ws.Cells.Item(67, 8).Value = 3228.5
$ws.Cells.Item(67, 10).Value = 3228.5
$ws.Cells.Item(67, 12).Value = 9685.5
$ws.Cells.Item(67, 14).Value = -11557.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(82, 8).Value = 30007
$ws.Cells.Item(82, 10).Value = 30007
$ws.Cells.Item(82, 12).Value = 30007
$ws.Cells.Item(82, 14).Value = -30773

$ws.Cells.Item(85, 8).Value = 30007
$ws.Cells.Item(85, 10).Value = 30007
$ws.Cells.Item(85, 12).Value = 30007
$ws.Cells.Item(85, 14).Value = -32659

$ws.Cells.Item(92, 8).Value = 5000
$ws.Cells.Item(92, 10).Value = 5000
$ws.Cells.Item(92, 12).Value = 5000
$ws.Cells.Item(92, 14).Value = -8744

$ws.Cells.Item(94, 8).Value = 0
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 12).Value = 0
$ws.Cells.Item(94, 14).ClearContents()

$ws.Cells.Item(122, 8).Value = 2721.4814
$ws.Cells.Item(122, 9).Value = 2743.9048
$ws.Cells.Item(122, 10).Value = 2643
$ws.Cells.Item(122, 11).Value = 8231.714399999999
$ws.Cells.Item(122, 12).Value = 7929
$ws.Cells.Item(122, 13).Value = -5781.714399999999
$ws.Cells.Item(122, 14).Value = -12829

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 17333.334
$ws.Cells.Item(7, 9).Value = 19700
$ws.Cells.Item(7, 10).Value = 5500
$ws.Cells.Item(7, 11).Value = 19700
$ws.Cells.Item(7, 12).Value = 5500
$ws.Cells.Item(7, 13).Value = -19588
$ws.Cells.Item(7, 14).Value = -5724

$ws.Cells.Item(16, 8).Value = 4625.125
$ws.Cells.Item(16, 9).Value = 3249.75
$ws.Cells.Item(16, 11).Value = 3249.75
$ws.Cells.Item(16, 13).Value = -3079.75

$ws.Cells.Item(40, 8).Value = 4159.3335
$ws.Cells.Item(40, 9).Value = 2885.8
$ws.Cells.Item(40, 10).Value = 5751.25
$ws.Cells.Item(40, 11).Value = 2885.8
$ws.Cells.Item(40, 12).Value = 5751.25
$ws.Cells.Item(40, 13).Value = -2749.8
$ws.Cells.Item(40, 14).Value = -6023.25

$ws.Cells.Item(122, 8).Value = 13086
$ws.Cells.Item(122, 9).Value = 14882.889
$ws.Cells.Item(122, 11).Value = 44648.667
$ws.Cells.Item(122, 13).Value = -42198.667

$ws.Cells.Item(126, 8).Value = 17333.334
$ws.Cells.Item(126, 9).Value = 19700
$ws.Cells.Item(126, 10).Value = 5500
$ws.Cells.Item(126, 11).Value = 59100
$ws.Cells.Item(126, 12).Value = 16500
$ws.Cells.Item(126, 13).Value = -56630
$ws.Cells.Item(126, 14).Value = -21440

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 31250696
$ws.Cells.Item(107, 9).Value = 726
$ws.Cells.Item(107, 11).Value = 2178
$ws.Cells.Item(107, 13).Value = -258
